# Expand import validation in MCC
# Applies the changes described by the commit to MCC_Data_Template.xlsx:
#  - adds an example/sample row (row 2) to the "animal data" sheet
#  - restricts the dropdown data validations to start at row 2 (leaving the
#    sample row itself out of the "apply to entire column" validation)
#  - relabels the "K" column header string and tidies a couple of leftover
#    cell styles
#  - updates the current selection / view state

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("instructions")
$wsData = $wb.Worksheets.Item("animal data")

# ---------------------------------------------------------------------------
# instructions sheet: drop the stray fill-only style that had been left on
# B17 / B18 (both just contain the "Dropdown" text with no real formatting)
# ---------------------------------------------------------------------------
$wsInstructions.Range("B17").Style = "Normal"
$wsInstructions.Range("B18").Style = "Normal"

# ---------------------------------------------------------------------------
# animal data sheet: fill in the sample/example row so users can see what
# valid data looks like (do this before the K1 header edit below so the new
# shared strings land in the same order the original author typed them in)
# ---------------------------------------------------------------------------
$wsData.Range("A2").Value = 1
$wsData.Range("B2").Value = 1
$wsData.Range("C2").Value = 1
$wsData.Range("D2").Value = 1
$wsData.Range("E2").Value = "1 - female"
$wsData.Range("F2").Value = 1
$wsData.Range("G2").Value = 1
$wsData.Range("H2").Value = 1
$wsData.Range("I2").Value = 1
$wsData.Range("J2").Value = "1 - assigned to U24 breeding colong"
$wsData.Range("K2").Value = "1 - available for transfer"
$wsData.Range("L2").Value = "1 - natal family group"
$wsData.Range("M2").Value = "1 - sibling experience only"
$wsData.Range("N2").Value = "1 - mated no offspring produced"
$wsData.Range("O2").Value = "1 - animal assigned to invasive study"

$wsData.Rows.Item(2).RowHeight = 45

# ---------------------------------------------------------------------------
# animal data sheet: header correction
# ---------------------------------------------------------------------------
$wsData.Range("K1").Value = "available to transfer"

# ---------------------------------------------------------------------------
# animal data sheet: the dropdown validations should no longer cover the
# header/example rows - restrict them to start at row 2 (data entry rows).
# Delete + re-Add (scoped to the new range) so the sqref in the saved XML
# actually changes; re-apply each validation in the same column order as
# before so the rules stay laid out the same way.
# ---------------------------------------------------------------------------
$wsData.Range("E1:E1048576").Validation.Delete()
$wsData.Range("E2:E1048576").Validation.Add(3, 1, 1, '"0 - male, 1 - female"')
$wsData.Range("E2:E1048576").Validation.IgnoreBlank = 1
$wsData.Range("E2:E1048576").Validation.InCellDropdown = 1
$wsData.Range("E2:E1048576").Validation.ShowInput = 1
$wsData.Range("E2:E1048576").Validation.ShowError = 1

$wsData.Range("J1:J1048576").Validation.Delete()
$wsData.Range("J2:J1048576").Validation.Add(3, 1, 1, '"0 - not assigned to U24 breeding colony, 1 - assigned to U24 breeding colong"')
$wsData.Range("J2:J1048576").Validation.IgnoreBlank = 1
$wsData.Range("J2:J1048576").Validation.InCellDropdown = 1
$wsData.Range("J2:J1048576").Validation.ShowInput = 1
$wsData.Range("J2:J1048576").Validation.ShowError = 1

$wsData.Range("K1:K1048576").Validation.Delete()
$wsData.Range("K2:K1048576").Validation.Add(3, 1, 1, '"0 - not available for transfer, 1 - available for transfer"')
$wsData.Range("K2:K1048576").Validation.IgnoreBlank = 1
$wsData.Range("K2:K1048576").Validation.InCellDropdown = 1
$wsData.Range("K2:K1048576").Validation.ShowInput = 1
$wsData.Range("K2:K1048576").Validation.ShowError = 1

$wsData.Range("L1:L1048576").Validation.Delete()
$wsData.Range("L2:L1048576").Validation.Add(3, 1, 1, '"0 - singly housed, 1 - natal family group, 2 - active breeding, 3 - social non breeding"')
$wsData.Range("L2:L1048576").Validation.IgnoreBlank = 1
$wsData.Range("L2:L1048576").Validation.InCellDropdown = 1
$wsData.Range("L2:L1048576").Validation.ShowInput = 1
$wsData.Range("L2:L1048576").Validation.ShowError = 1

$wsData.Range("M1:M1048576").Validation.Delete()
$wsData.Range("M2:M1048576").Validation.Add(3, 1, 1, '"0 - no experience, 1 - sibling experience only, 2 - non successful offspring, 3 - successful rearing of offspring"')
$wsData.Range("M2:M1048576").Validation.IgnoreBlank = 1
$wsData.Range("M2:M1048576").Validation.InCellDropdown = 1
$wsData.Range("M2:M1048576").Validation.ShowInput = 1
$wsData.Range("M2:M1048576").Validation.ShowError = 1

$wsData.Range("N1:N1048576").Validation.Delete()
$wsData.Range("N2:N1048576").Validation.Add(3, 1, 1, '"0 - no mating opportunity, 1 - mated no offspring produced, 2 - successful offspring produced, 3 - hormonal birth control, 4 - sterilized"')
$wsData.Range("N2:N1048576").Validation.IgnoreBlank = 1
$wsData.Range("N2:N1048576").Validation.InCellDropdown = 1
$wsData.Range("N2:N1048576").Validation.ShowInput = 1
$wsData.Range("N2:N1048576").Validation.ShowError = 1

$wsData.Range("O1:O1048576").Validation.Delete()
$wsData.Range("O2:O1048576").Validation.Add(3, 1, 1, '"0 - naive animal, 1 - animal assigned to invasive study"')
$wsData.Range("O2:O1048576").Validation.IgnoreBlank = 1
$wsData.Range("O2:O1048576").Validation.InCellDropdown = 1
$wsData.Range("O2:O1048576").Validation.ShowInput = 1
$wsData.Range("O2:O1048576").Validation.ShowError = 1

# ---------------------------------------------------------------------------
# page setup / view state
# ---------------------------------------------------------------------------
$wsData.PageSetup.Orientation = 1

$wsData.Range("L15").Select()

Write-Host "Done"
